$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title update
$ws.Range("A1").Value = 'Datos actualizados a 1 de Abril de 2020 a las 14:50'

# Row 8
$ws.Cells.Item(8, 2).Value = 73217
$ws.Cells.Item(8, 3).Value = 1409
$ws.Cells.Item(8, 5).Value = 56315
$ws.Cells.Item(8, 7).Value = 27
$ws.Cells.Item(8, 8).Value = 802

# Row 24
$ws.Cells.Item(24, 2).Value = 4798
$ws.Cells.Item(24, 3).Value = 157
$ws.Cells.Item(24, 5).Value = 4742
$ws.Cells.Item(24, 7).Value = 4
$ws.Cells.Item(24, 8).Value = 43

# Row 32
$ws.Cells.Item(32, 2).Value = 2420
$ws.Cells.Item(32, 3).Value = 109
$ws.Cells.Item(32, 5).Value = 2377
$ws.Cells.Item(32, 7).Value = 3
$ws.Cells.Item(32, 8).Value = 36

# Row 37
$ws.Cells.Item(37, 2).Value = 2071
$ws.Cells.Item(37, 3).Value = 133
$ws.Cells.Item(37, 5).Value = 1963

# Row 39
$ws.Cells.Item(39, 1).Value = 'Arabia Saudita'
$ws.Cells.Item(39, 2).Value = 1720
$ws.Cells.Item(39, 3).Value = 157
$ws.Cells.Item(39, 4).Value = 264
$ws.Cells.Item(39, 5).Value = 1440
$ws.Cells.Item(39, 6).Value = 31
$ws.Cells.Item(39, 7).Value = 6
$ws.Cells.Item(39, 8).Value = 16

# Row 40
$ws.Cells.Item(40, 1).Value = 'Indonesia'
$ws.Cells.Item(40, 2).Value = 1677
$ws.Cells.Item(40, 3).Value = 149
$ws.Cells.Item(40, 4).Value = 103
$ws.Cells.Item(40, 5).Value = 1417
$ws.Cells.Item(40, 7).Value = 21
$ws.Cells.Item(40, 8).Value = 157

# Row 41
$ws.Cells.Item(41, 1).Value = 'India'
$ws.Cells.Item(41, 2).Value = 1637
$ws.Cells.Item(41, 3).Value = 240
$ws.Cells.Item(41, 4).Value = 148
$ws.Cells.Item(41, 5).Value = 1444
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(41, 7).Value = 10
$ws.Cells.Item(41, 8).Value = 45

# Row 117
$ws.Cells.Item(117, 4).Value = 3
$ws.Cells.Item(117, 5).Value = 97
$ws.Cells.Item(117, 7).Value = 1
$ws.Cells.Item(117, 8).Value = 9

# Row 121
$ws.Cells.Item(121, 1).Value = 'Kenia'
$ws.Cells.Item(121, 2).Value = 81
$ws.Cells.Item(121, 3).Value = 22
$ws.Cells.Item(121, 4).Value = 3
$ws.Cells.Item(121, 5).Value = 77
$ws.Cells.Item(121, 6).Value = 2
$ws.Cells.Item(121, 8).Value = 1

# Row 122
$ws.Cells.Item(122, 1).Value = 'Ruanda'
$ws.Cells.Item(122, 2).Value = 75
$ws.Cells.Item(122, 3).Value = 0
$ws.Cells.Item(122, 4).Value = 0
$ws.Cells.Item(122, 5).Value = 75
$ws.Cells.Item(122, 6).Value = 0
$ws.Cells.Item(122, 8).Value = 0

# Row 123
$ws.Cells.Item(123, 1).Value = 'Paraguay'
$ws.Cells.Item(123, 3).Value = 4
$ws.Cells.Item(123, 4).Value = 1
$ws.Cells.Item(123, 5).Value = 65
$ws.Cells.Item(123, 6).Value = 3
$ws.Cells.Item(123, 8).Value = 3

# Row 124
$ws.Cells.Item(124, 1).Value = 'Gibraltar'
$ws.Cells.Item(124, 2).Value = 69
$ws.Cells.Item(124, 4).Value = 34
$ws.Cells.Item(124, 5).Value = 35

# Row 125
$ws.Cells.Item(125, 1).Value = 'Liechtenstein'
$ws.Cells.Item(125, 2).Value = 68
$ws.Cells.Item(125, 3).Value = 0
$ws.Cells.Item(125, 5).Value = 68

# Row 126
$ws.Cells.Item(126, 1).Value = 'Isla de Man'
$ws.Cells.Item(126, 2).Value = 65
$ws.Cells.Item(126, 3).Value = 5
$ws.Cells.Item(126, 4).Value = 0
$ws.Cells.Item(126, 5).Value = 65
$ws.Cells.Item(126, 6).Value = 0
$ws.Cells.Item(126, 8).Value = 0

# Row 142
$ws.Cells.Item(142, 1).Value = 'Guam'
$ws.Cells.Item(142, 6).Value = 0

# Row 143
$ws.Cells.Item(143, 1).Value = 'El Salvador'
$ws.Cells.Item(143, 6).Value = 4

# Row 146
$ws.Cells.Item(146, 1).Value = 'Republica de Yibuti'
$ws.Cells.Item(146, 3).Value = 0

# Row 147
$ws.Cells.Item(147, 1).Value = 'Guinea'
$ws.Cells.Item(147, 3).Value = 8

# Row 158
$ws.Cells.Item(158, 1).Value = 'Birmania'
$ws.Cells.Item(158, 4).Value = 0
$ws.Cells.Item(158, 8).Value = 1

# Row 159
$ws.Cells.Item(159, 1).Value = 'Guinea Ecuatorial'
$ws.Cells.Item(159, 4).Value = 1
$ws.Cells.Item(159, 8).Value = 0

# Row 164
$ws.Cells.Item(164, 1).Value = 'Namibia'
$ws.Cells.Item(164, 2).Value = 14
$ws.Cells.Item(164, 3).Value = 3
$ws.Cells.Item(164, 4).Value = 2

# Row 165
$ws.Cells.Item(165, 1).Value = 'Santa Lucia'
$ws.Cells.Item(165, 2).Value = 13
$ws.Cells.Item(165, 4).Value = 1

# Row 166
$ws.Cells.Item(166, 1).Value = 'Dominica'
$ws.Cells.Item(166, 5).Value = 12
$ws.Cells.Item(166, 8).Value = 0

# Row 167
$ws.Cells.Item(167, 1).Value = 'Guyana'
$ws.Cells.Item(167, 2).Value = 12
$ws.Cells.Item(167, 4).Value = 0
$ws.Cells.Item(167, 5).Value = 10
$ws.Cells.Item(167, 8).Value = 2

# Row 169
$ws.Cells.Item(169, 1).Value = 'Laos'

# Row 171
$ws.Cells.Item(171, 1).Value = 'Surinam'

# Row 178
$ws.Cells.Item(178, 1).Value = 'Guinea-Bisau'

# Row 179
$ws.Cells.Item(179, 1).Value = 'Mozambique'

# Row 180
$ws.Cells.Item(180, 1).Value = 'San Cristobal y Nieves'

# Row 186
$ws.Cells.Item(186, 1).Value = 'Santa Sede'
$ws.Cells.Item(186, 3).Value = 0

# Row 187
$ws.Cells.Item(187, 1).Value = 'Liberia'
$ws.Cells.Item(187, 3).Value = 3

# Row 188
$ws.Cells.Item(188, 1).Value = 'Cabo Verde'
$ws.Cells.Item(188, 4).Value = 0
$ws.Cells.Item(188, 8).Value = 1

# Row 189
$ws.Cells.Item(189, 1).Value = 'San Bartolome'
$ws.Cells.Item(189, 4).Value = 1
$ws.Cells.Item(189, 8).Value = 0

# Row 191
$ws.Cells.Item(191, 1).Value = 'Fiyi'

# Row 192
$ws.Cells.Item(192, 1).Value = 'Montserrat'
